# Scheduled Leve-profit refresh: push updated currentAveragePrice / leve-profit
# figures (columns H-N) across all eight crafting-sheet tables. Only numeric
# market/profit cells are touched; leve metadata (A-G) is left untouched.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 20: Shut Up and Take My Gil - Ash Wand
$ws.Range("H20").Value = 4996.6665
$ws.Range("I20").Value = 4996.6665
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 4996.6665
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -4766.6665

# Row 28: The Writing Is Not on the Wall - Enchanted Silver Ink
$ws.Range("H28").Value = 1999.75
$ws.Range("I28").Value = 1999.6666
$ws.Range("J28").Value = 2000
$ws.Range("K28").Value = 1999.6666
$ws.Range("L28").Value = 2000
$ws.Range("M28").Value = -1514.6666
$ws.Range("N28").Value = -2970

# Row 35: Conspicuous Conjuration - Whispering Ash Wand
$ws.Range("H35").Value = 4996.6665
$ws.Range("I35").Value = 4996.6665
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 4996.6665
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -4617.6665

# Row 63: Summoning for Dummies - Archaeoskin Codex
$ws.Range("H63").Value = 35000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 35000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 35000
$ws.Range("N63").Value = -36248

# Row 66: Summoning the Courage to Be Different (L) - Archaeoskin Codex
$ws.Range("H66").Value = 35000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 35000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 105000
$ws.Range("N66").Value = -111240

# Row 96: Scroll Down - Grade 1 Reisui of Intelligence
$ws.Range("H96").Value = 446.33334
$ws.Range("I96").Value = 533.3333
$ws.Range("J96").Value = 359.33334
$ws.Range("K96").Value = 1599.9999
$ws.Range("L96").Value = 1078.00002
$ws.Range("M96").Value = -226.9999
$ws.Range("N96").Value = -3824.00002

# Row 115: 5-bell Energy - Competent Craftsman's Syrup
$ws.Range("H115").Value = 1085
$ws.Range("I115").Value = 1085
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 3255
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -1688

# Row 116: Growing Up - Growth Formula Kappa
$ws.Range("H116").Value = 4753
$ws.Range("I116").Value = 4930
$ws.Range("J116").Value = 4399
$ws.Range("K116").Value = 4930
$ws.Range("L116").Value = 4399
$ws.Range("M116").Value = -1488
$ws.Range("N116").Value = -11283

# Row 134: Binding Spells - Crocodileskin Index
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# Row 137: Cutting Edge of Culinary Quality - Magnesia Whetstone
$ws.Range("H137").Value = 8179.8096
$ws.Range("I137").Value = 5085.846
$ws.Range("J137").Value = 13207.5
$ws.Range("K137").Value = 15257.538
$ws.Range("L137").Value = 39622.5
$ws.Range("M137").Value = -12707.538
$ws.Range("N137").Value = -44722.5

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 19: Stadium Envy - Bronze Gauntlets
$ws.Range("H19").Value = 1499
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1499
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 1499
$ws.Range("N19").Value = -1957

# Row 32: Ingot We Trust - Steel Ingot
$ws.Range("H32").Value = 2000
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1713

# Row 63: Rivets Run through It - Mythrite Rivets
$ws.Range("H63").Value = 7950
$ws.Range("I63").Value = 7950
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 7950
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -7264

# Row 66: A Riveting Revival (L) - Mythrite Rivets
$ws.Range("H66").Value = 7950
$ws.Range("I66").Value = 7950
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 39750
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -36318

# Row 122: Haste for High Durium - High Durium Nugget
$ws.Range("H122").Value = 4998.3335
$ws.Range("I122").Value = 4999
$ws.Range("J122").Value = 4997
$ws.Range("K122").Value = 14997
$ws.Range("L122").Value = 14991
$ws.Range("M122").Value = -12547
$ws.Range("N122").Value = -19891

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 82: Spirituality Inspector - Titanium Lump Hammer
$ws.Range("H82").Value = 9613.857
$ws.Range("I82").Value = 2382.8333
$ws.Range("J82").Value = 53000
$ws.Range("K82").Value = 2382.8333
$ws.Range("L82").Value = 53000
$ws.Range("M82").Value = -1999.8333
$ws.Range("N82").Value = -53766

# Row 85: The Clamor for Hammers (L) - Titanium Lump Hammer
$ws.Range("H85").Value = 9613.857
$ws.Range("I85").Value = 2382.8333
$ws.Range("J85").Value = 53000
$ws.Range("K85").Value = 2382.8333
$ws.Range("L85").Value = 53000
$ws.Range("M85").Value = -1056.8333
$ws.Range("N85").Value = -55652

# Row 86: Through Thick and Thin - Adamantite Nugget
$ws.Range("H86").Value = 2498.3333
$ws.Range("I86").Value = 2498.3333
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2498.3333
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1375.3333

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) - Adamantite Nugget
$ws.Range("H89").Value = 2498.3333
$ws.Range("I89").Value = 2498.3333
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 12491.6665
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -6875.666499999999

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 2: In with the New - Bone Harpoon
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# Row 31: Wall Not Found - Walnut Lumber
$ws.Range("H31").Value = 6285.7144
$ws.Range("I31").Value = 5000
$ws.Range("J31").Value = 8000
$ws.Range("K31").Value = 5000
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = -4705
$ws.Range("N31").Value = -8590

# Row 32: Daddy's Little Girl - Viper-crested Round Shield
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()

# Row 34: Armoires of the Rich and Famous - Walnut Lumber
$ws.Range("H34").Value = 6285.7144
$ws.Range("I34").Value = 5000
$ws.Range("J34").Value = 8000
$ws.Range("K34").Value = 5000
$ws.Range("L34").Value = 8000
$ws.Range("M34").Value = -4798
$ws.Range("N34").Value = -8404

# Row 35: Storm of Swords - Elm Macuahuitl
$ws.Range("H35").Value = 1354
$ws.Range("I35").Value = 1354
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1354
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1060

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 7: It's Always Sunny in Vylbrand - Raisins
$ws.Range("H7").Value = 13
$ws.Range("I7").Value = 13
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 39
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 73
$ws.Range("N7").ClearContents()

# Row 103: West Meats East - Nomad Meat Pie
$ws.Range("H103").Value = 358.33334
$ws.Range("I103").Value = 87.5
$ws.Range("J103").Value = 900
$ws.Range("K103").Value = 262.5
$ws.Range("L103").Value = 2700
$ws.Range("M103").Value = 616.5
$ws.Range("N103").Value = -4458

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 46: Burning the Midnight Oil - Fire Brand
$ws.Range("H46").Value = 5000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5312
$ws.Range("M46").ClearContents()

# Row 92: Play It by Ear - Triphane Earrings of Healing
$ws.Range("H92").Value = 2050.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 2050.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 2050.5
$ws.Range("N92").Value = -5794.5

# Row 102: Put the Metal to the Peddle - Durium Ingot
$ws.Range("H102").Value = 3332.8333
$ws.Range("I102").Value = 2765.6667
$ws.Range("J102").Value = 3900
$ws.Range("K102").Value = 2765.6667
$ws.Range("L102").Value = 3900
$ws.Range("M102").Value = -1143.6667
$ws.Range("N102").Value = -7144

# Row 122: Awarding Academic Excellence - Ametrine
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 32: Men Who Scare Up Goats - Goatskin Targe
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -683

# Row 61: Spelling Me Softly - Raptor Leather
$ws.Range("H61").Value = 2966.3333
$ws.Range("I61").Value = 2799.5
$ws.Range("J61").Value = 3300
$ws.Range("K61").Value = 2799.5
$ws.Range("L61").Value = 3300
$ws.Range("M61").Value = -2597.5
$ws.Range("N61").Value = -3704

# Row 100: Tiger in the Sack - Tiger Leather
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()

# Row 113: Peace in Rest - Atrociraptor Leather
$ws.Range("H113").Value = 2966.3333
$ws.Range("I113").Value = 2799.5
$ws.Range("J113").Value = 3300
$ws.Range("K113").Value = 2799.5
$ws.Range("L113").Value = 3300
$ws.Range("M113").Value = -629.5
$ws.Range("N113").Value = -7640

# Row 132: Tenets of Tanning - Silver Lobo Leather
$ws.Range("H132").Value = 12773.454
$ws.Range("I132").Value = 8438.5
$ws.Range("J132").Value = 24333.334
$ws.Range("K132").Value = 25315.5
$ws.Range("L132").Value = 73000.00199999999
$ws.Range("M132").Value = -22785.5
$ws.Range("N132").Value = -78060.00199999999

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 92: Modest Beginnings - Bloodhempen Culottes of Casting
$ws.Range("H92").Value = 34000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 34000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 34000
$ws.Range("N92").Value = -38992

# Row 101: Who War It Better - Serge Hose of Aiming
$ws.Range("H101").Value = 18049.75
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 18049.75
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 18049.75
$ws.Range("N101").Value = -24539.75

# Row 107: Flax Wax - Bright Linen Yarn
$ws.Range("H107").Value = 998.8
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 998.8
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2996.4
$ws.Range("N107").Value = -6836.4

# Row 111: Legs for Days - Iridescent Bottoms of Scouting
$ws.Range("H111").Value = 38650
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 38650
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 38650
$ws.Range("N111").Value = -46830
